$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.06832672688341
$ws.Cells.Item(2, 4).Value = 1.07413319376261
$ws.Cells.Item(2, 5).Value = 1.062824050288019
$ws.Cells.Item(2, 6).Value = 1.084164820100711
$ws.Cells.Item(2, 9).Value = 1.061370468564422
$ws.Cells.Item(2, 10).Value = 1.073267100808933
$ws.Cells.Item(2, 11).Value = 1.076823213345874
$ws.Cells.Item(2, 12).Value = 1.065544351886836
$ws.Cells.Item(2, 13).Value = 1.086828549394838
$ws.Cells.Item(2, 14).Value = 1.074791263197538
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.06969516042007
$ws.Cells.Item(3, 4).Value = 1.075274563170288
$ws.Cells.Item(3, 5).Value = 1.064030264119763
$ws.Cells.Item(3, 6).Value = 1.085454180209984
$ws.Cells.Item(3, 9).Value = 1.061912305286399
$ws.Cells.Item(3, 10).Value = 1.074290282042476
$ws.Cells.Item(3, 11).Value = 1.077780853223571
$ws.Cells.Item(3, 12).Value = 1.066564473865797
$ws.Cells.Item(3, 13).Value = 1.087935730445933
$ws.Cells.Item(3, 14).Value = 1.075815897465794
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.070579106151959
$ws.Cells.Item(4, 4).Value = 1.076011830222849
$ws.Cells.Item(4, 5).Value = 1.0648089010467
$ws.Cells.Item(4, 6).Value = 1.08628756955649
$ws.Cells.Item(4, 9).Value = 1.062260890813131
$ws.Cells.Item(4, 10).Value = 1.074950333521312
$ws.Cells.Item(4, 11).Value = 1.07839864941301
$ws.Cells.Item(4, 12).Value = 1.067222129601094
$ws.Cells.Item(4, 13).Value = 1.088650678396311
$ws.Cells.Item(4, 14).Value = 1.076476886293439
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.070950358566146
$ws.Cells.Item(5, 4).Value = 1.076321476317795
$ws.Cells.Item(5, 5).Value = 1.065135799059161
$ws.Cells.Item(5, 6).Value = 1.086637712594998
$ws.Cells.Item(5, 9).Value = 1.062406955479014
$ws.Cells.Item(5, 10).Value = 1.075227341339786
$ws.Cells.Item(5, 11).Value = 1.078657929560317
$ws.Cells.Item(5, 12).Value = 1.067498030941228
$ws.Cells.Item(5, 13).Value = 1.088950894437825
$ws.Cells.Item(5, 14).Value = 1.076754287494789
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.071012672640649
$ws.Cells.Item(6, 4).Value = 1.076373449735742
$ws.Cells.Item(6, 5).Value = 1.06519066097015
$ws.Cells.Item(6, 6).Value = 1.086696490748074
$ws.Cells.Item(6, 9).Value = 1.062431452268817
$ws.Cells.Item(6, 10).Value = 1.07527382428539
$ws.Cells.Item(6, 11).Value = 1.078701438075036
$ws.Cells.Item(6, 12).Value = 1.067544322290431
$ws.Cells.Item(6, 13).Value = 1.089001281747551
$ws.Cells.Item(6, 14).Value = 1.076800836451506
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.070584068245514
$ws.Cells.Item(7, 4).Value = 1.076015968906378
$ws.Cells.Item(7, 5).Value = 1.064813270799981
$ws.Cells.Item(7, 6).Value = 1.086292249019332
$ws.Cells.Item(7, 9).Value = 1.062262844421813
$ws.Cells.Item(7, 10).Value = 1.074954036785131
$ws.Cells.Item(7, 11).Value = 1.078402115656208
$ws.Cells.Item(7, 12).Value = 1.06722581846989
$ws.Cells.Item(7, 13).Value = 1.088654691260029
$ws.Cells.Item(7, 14).Value = 1.076480594816318
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.068789514258682
$ws.Cells.Item(8, 4).Value = 1.074519191565571
$ws.Cells.Item(8, 5).Value = 1.063232085075781
$ws.Cells.Item(8, 6).Value = 1.084600757235982
$ws.Cells.Item(8, 9).Value = 1.061554004498714
$ws.Cells.Item(8, 10).Value = 1.073613309744933
$ws.Cells.Item(8, 11).Value = 1.077147240013553
$ws.Cells.Item(8, 12).Value = 1.06588961331032
$ws.Cells.Item(8, 13).Value = 1.087203033890502
$ws.Cells.Item(8, 14).Value = 1.075137963789939
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.065615349611882
$ws.Cells.Item(9, 4).Value = 1.071871716686891
$ws.Cells.Item(9, 5).Value = 1.060431318341849
$ws.Cells.Item(9, 6).Value = 1.081612929780096
$ws.Cells.Item(9, 9).Value = 1.060289353834373
$ws.Cells.Item(9, 10).Value = 1.071235135659737
$ws.Cells.Item(9, 11).Value = 1.074921551560948
$ws.Cells.Item(9, 12).Value = 1.063516204943581
$ws.Cells.Item(9, 13).Value = 1.084633555705241
$ws.Cells.Item(9, 14).Value = 1.072756412424761
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.063490834924698
$ws.Cells.Item(10, 4).Value = 1.07009976069398
$ws.Cells.Item(10, 5).Value = 1.05855404505182
$ws.Cells.Item(10, 6).Value = 1.079615904000094
$ws.Cells.Item(10, 9).Value = 1.059435605544277
$ws.Cells.Item(10, 10).Value = 1.069638877170021
$ws.Cells.Item(10, 11).Value = 1.073427791218415
$ws.Cells.Item(10, 12).Value = 1.061920949890704
$ws.Cells.Item(10, 13).Value = 1.08291257919462
$ws.Cells.Item(10, 14).Value = 1.071157887064896
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.062568817225964
$ws.Cells.Item(11, 4).Value = 1.069330768459502
$ws.Cells.Item(11, 5).Value = 1.057738695053263
$ws.Cells.Item(11, 6).Value = 1.07874988339727
$ws.Cells.Item(11, 9).Value = 1.059063359330308
$ws.Cells.Item(11, 10).Value = 1.068945050153974
$ws.Cells.Item(11, 11).Value = 1.07277855345181
$ws.Cells.Item(11, 12).Value = 1.061227037251216
$ws.Cells.Item(11, 13).Value = 1.082165419527929
$ws.Cells.Item(11, 14).Value = 1.070463074734904
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.062226017494912
$ws.Cells.Item(12, 4).Value = 1.069044866129703
$ws.Cells.Item(12, 5).Value = 1.057435458836371
$ws.Cells.Item(12, 6).Value = 1.078428004318291
$ws.Cells.Item(12, 9).Value = 1.058924701398635
$ws.Cells.Item(12, 10).Value = 1.068686930045567
$ws.Cells.Item(12, 11).Value = 1.072537027454938
$ws.Cells.Item(12, 12).Value = 1.060968806802419
$ws.Cells.Item(12, 13).Value = 1.081887590396373
$ws.Cells.Item(12, 14).Value = 1.070204588066335
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.062299563861321
$ws.Cells.Item(13, 4).Value = 1.069106205200999
$ws.Cells.Item(13, 5).Value = 1.057500521287689
$ws.Cells.Item(13, 6).Value = 1.078497057639665
$ws.Cells.Item(13, 9).Value = 1.058954461658163
$ws.Cells.Item(13, 10).Value = 1.068742315962884
$ws.Cells.Item(13, 11).Value = 1.072588852435822
$ws.Cells.Item(13, 12).Value = 1.061024219948384
$ws.Cells.Item(13, 13).Value = 1.081947199386724
$ws.Cells.Item(13, 14).Value = 1.070260052638008
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.062540487891253
$ws.Cells.Item(14, 4).Value = 1.069307141108415
$ws.Cells.Item(14, 5).Value = 1.057713637221203
$ws.Cells.Item(14, 6).Value = 1.078723280895775
$ws.Cells.Item(14, 9).Value = 1.059051905779772
$ws.Cells.Item(14, 10).Value = 1.068923722096558
$ws.Cells.Item(14, 11).Value = 1.072758596419117
$ws.Cells.Item(14, 12).Value = 1.061205701677423
$ws.Cells.Item(14, 13).Value = 1.082142460252544
$ws.Cells.Item(14, 14).Value = 1.070441716389199
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.062688886336824
$ws.Cells.Item(15, 4).Value = 1.069430909157547
$ws.Cells.Item(15, 5).Value = 1.057844894572497
$ws.Cells.Item(15, 6).Value = 1.078862637769957
$ws.Cells.Item(15, 9).Value = 1.059111892656756
$ws.Cells.Item(15, 10).Value = 1.069035438979369
$ws.Cells.Item(15, 11).Value = 1.072863132081899
$ws.Cells.Item(15, 12).Value = 1.061317454718444
$ws.Cells.Item(15, 13).Value = 1.082262726897828
$ws.Cells.Item(15, 14).Value = 1.070553591922799
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.063551981556181
$ws.Cells.Item(16, 4).Value = 1.070150759399784
$ws.Cells.Item(16, 5).Value = 1.058608104369983
$ws.Cells.Item(16, 6).Value = 1.079673351145706
$ws.Cells.Item(16, 9).Value = 1.0594602559171
$ws.Cells.Item(16, 10).Value = 1.069684868145013
$ws.Cells.Item(16, 11).Value = 1.073470827361394
$ws.Cells.Item(16, 12).Value = 1.061966935556323
$ws.Cells.Item(16, 13).Value = 1.082962123844365
$ws.Cells.Item(16, 14).Value = 1.071203943352348
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.064092813881038
$ws.Cells.Item(17, 4).Value = 1.0706018371665
$ws.Cells.Item(17, 5).Value = 1.059086177568338
$ws.Cells.Item(17, 6).Value = 1.080181539319134
$ws.Cells.Item(17, 9).Value = 1.059678085206846
$ws.Cells.Item(17, 10).Value = 1.070091527975664
$ws.Cells.Item(17, 11).Value = 1.073851364435683
$ws.Cells.Item(17, 12).Value = 1.062373488271175
$ws.Cells.Item(17, 13).Value = 1.083400307004811
$ws.Cells.Item(17, 14).Value = 1.071611180686604
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.064408071193981
$ws.Cells.Item(18, 4).Value = 1.070864777263379
$ws.Cells.Item(18, 5).Value = 1.059364790907147
$ws.Cells.Item(18, 6).Value = 1.080477832308379
$ws.Cells.Item(18, 9).Value = 1.059804893841818
$ws.Cells.Item(18, 10).Value = 1.070328471776367
$ws.Cells.Item(18, 11).Value = 1.074073091194475
$ws.Cells.Item(18, 12).Value = 1.062610319435804
$ws.Cells.Item(18, 13).Value = 1.083655702560913
$ws.Cells.Item(18, 14).Value = 1.071848460974681
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.064515531917523
$ws.Cells.Item(19, 4).Value = 1.070954405059494
$ws.Cells.Item(19, 5).Value = 1.059459750571604
$ws.Cells.Item(19, 6).Value = 1.080578839607381
$ws.Cells.Item(19, 9).Value = 1.059848090429546
$ws.Cells.Item(19, 10).Value = 1.070409220579745
$ws.Cells.Item(19, 11).Value = 1.074148654740515
$ws.Cells.Item(19, 12).Value = 1.062691021337578
$ws.Cells.Item(19, 13).Value = 1.083742753907061
$ws.Cells.Item(19, 14).Value = 1.071929324450622
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.064034808540499
$ws.Cells.Item(20, 4).Value = 1.070553457994904
$ws.Cells.Item(20, 5).Value = 1.059034909574521
$ws.Cells.Item(20, 6).Value = 1.080127028466834
$ws.Cells.Item(20, 9).Value = 1.059654739811286
$ws.Cells.Item(20, 10).Value = 1.070047923539631
$ws.Cells.Item(20, 11).Value = 1.07381056065934
$ws.Cells.Item(20, 12).Value = 1.062329900518559
$ws.Cells.Item(20, 13).Value = 1.083353313716448
$ws.Cells.Item(20, 14).Value = 1.07156751432727
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.06246955069926
$ws.Cells.Item(21, 4).Value = 1.069247977833675
$ws.Cells.Item(21, 5).Value = 1.05765089040388
$ws.Cells.Item(21, 6).Value = 1.078656669366559
$ws.Cells.Item(21, 9).Value = 1.059023221683274
$ws.Cells.Item(21, 10).Value = 1.068870313651358
$ws.Cells.Item(21, 11).Value = 1.072708621307217
$ws.Cells.Item(21, 12).Value = 1.06115227314242
$ws.Cells.Item(21, 13).Value = 1.082084969141027
$ws.Cells.Item(21, 14).Value = 1.070388232097881
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.061483548255268
$ws.Cells.Item(22, 4).Value = 1.068425638336593
$ws.Cells.Item(22, 5).Value = 1.056778506985324
$ws.Cells.Item(22, 6).Value = 1.077731033460327
$ws.Cells.Item(22, 9).Value = 1.05862390819647
$ws.Cells.Item(22, 10).Value = 1.068127575649094
$ws.Cells.Item(22, 11).Value = 1.072013643957522
$ws.Cells.Item(22, 12).Value = 1.060409069563871
$ws.Cells.Item(22, 13).Value = 1.081285767698284
$ws.Cells.Item(22, 14).Value = 1.069644439322461
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.0620064258816
$ws.Cells.Item(23, 4).Value = 1.068861723128967
$ws.Cells.Item(23, 5).Value = 1.057241184274151
$ws.Cells.Item(23, 6).Value = 1.078221842756343
$ws.Cells.Item(23, 9).Value = 1.058835806516488
$ws.Cells.Item(23, 10).Value = 1.068521537707236
$ws.Cells.Item(23, 11).Value = 1.072382269588012
$ws.Cells.Item(23, 12).Value = 1.060803321653609
$ws.Cells.Item(23, 13).Value = 1.081709606545357
$ws.Cells.Item(23, 14).Value = 1.070038960851913
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.064061019277681
$ws.Cells.Item(24, 4).Value = 1.070575318972437
$ws.Cells.Item(24, 5).Value = 1.0590580761111
$ws.Cells.Item(24, 6).Value = 1.080151659959129
$ws.Cells.Item(24, 9).Value = 1.059665289355152
$ws.Cells.Item(24, 10).Value = 1.070067627291738
$ws.Cells.Item(24, 11).Value = 1.073828998852391
$ws.Cells.Item(24, 12).Value = 1.062349596887089
$ws.Cells.Item(24, 13).Value = 1.083374548545091
$ws.Cells.Item(24, 14).Value = 1.071587246060965
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.066437400451404
$ws.Cells.Item(25, 4).Value = 1.072557361119881
$ws.Cells.Item(25, 5).Value = 1.061157138296981
$ws.Cells.Item(25, 6).Value = 1.082386238917779
$ws.Cells.Item(25, 9).Value = 1.060618160568784
$ws.Cells.Item(25, 10).Value = 1.071851834279099
$ws.Cells.Item(25, 11).Value = 1.075498683043786
$ws.Cells.Item(25, 12).Value = 1.064132053130746
$ws.Cells.Item(25, 13).Value = 1.085299215408151
$ws.Cells.Item(25, 14).Value = 1.073373986826899

Write-Host "Updated vm_pu.xlsx values for 380 kV case"
